# Insert a new data row at row 421 (pushes old rows 421..518 down to 422..519,
# growing the sheet's used range from A1:R518 to A1:R519).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A421").EntireRow.Insert()

# Populate the newly inserted row 421 with the new record. This mirrors
# the row that used to sit at 421 (same market/region/category/variety/
# quality/unit/origin/kg-or-units/classification) but with its own date,
# volume, weighted average price and $/Kg.
$ws.Cells.Item(421, 1).Value = 6
$ws.Cells.Item(421, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(421, 3).Value = "Metropolitana"
$ws.Cells.Item(421, 4).Value = "2023-01-17"
$ws.Cells.Item(421, 5).Value = 13
$ws.Cells.Item(421, 6).Value = 100112032
$ws.Cells.Item(421, 7).Value = "Zapallo italiano"
$ws.Cells.Item(421, 8).Value = "Sin especificar"
$ws.Cells.Item(421, 9).Value = "Primera"
$ws.Cells.Item(421, 10).Value = 1250
$ws.Cells.Item(421, 11).Value = 9000
$ws.Cells.Item(421, 12).Value = 10000
$ws.Cells.Item(421, 13).Value = 9544
$ws.Cells.Item(421, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(421, 15).Value = "Región Metropolitana"
$ws.Cells.Item(421, 16).Value = 191
$ws.Cells.Item(421, 17).Value = 50
$ws.Cells.Item(421, 18).Value = "Hortaliza"
